$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()

$ws.Range("H112").Value = 6962.8125
$ws.Range("J112").Value = 6962.8125
$ws.Range("L112").Value = 20888.4375
$ws.Range("N112").Value = -23104.4375

$ws.Range("H116").Value = 9384.933999999999
$ws.Range("I116").Value = 11345.818
$ws.Range("K116").Value = 11345.818
$ws.Range("M116").Value = -7903.817999999999

$ws.Range("H127").Value = 1178.45
$ws.Range("I127").Value = 699.7143
$ws.Range("J127").Value = 1436.2307
$ws.Range("K127").Value = 2099.1429
$ws.Range("L127").Value = 4308.6921
$ws.Range("M127").Value = 2860.8571
$ws.Range("N127").Value = -14228.6921

$ws.Range("H129").Value = 1064.1637
$ws.Range("I129").Value = 551.55554
$ws.Range("J129").Value = 1164.4565
$ws.Range("K129").Value = 1654.66662
$ws.Range("L129").Value = 3493.3695
$ws.Range("M129").Value = 3345.33338
$ws.Range("N129").Value = -13493.3695

$ws.Range("H132").Value = 6246.7915
$ws.Range("I132").Value = 5125
$ws.Range("K132").Value = 15375
$ws.Range("M132").Value = -12845

$ws.Range("H138").Value = 2296.9592
$ws.Range("I138").Value = 2053.652
$ws.Range("J138").Value = 2512.1924
$ws.Range("K138").Value = 6160.956
$ws.Range("L138").Value = 7536.5772
$ws.Range("M138").Value = -1020.956
$ws.Range("N138").Value = -17816.5772

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 500
$ws.Range("I3").Value = 500
$ws.Range("K3").Value = 500
$ws.Range("M3").Value = -385

$ws.Range("H30").Value = 34084
$ws.Range("I30").Value = 5205
$ws.Range("J30").Value = 53336.668
$ws.Range("K30").Value = 5205
$ws.Range("L30").Value = 53336.668
$ws.Range("M30").Value = -5055
$ws.Range("N30").Value = -53636.668

$ws.Range("H61").Value = 2966.2258
$ws.Range("I61").Value = 2423.353
$ws.Range("J61").Value = 3625.4285
$ws.Range("K61").Value = 2423.353
$ws.Range("L61").Value = 3625.4285
$ws.Range("M61").Value = -2211.353
$ws.Range("N61").Value = -4049.4285

$ws.Range("H74").Value = 1843.9375
$ws.Range("J74").Value = 2988
$ws.Range("L74").Value = 2988
$ws.Range("N74").Value = -4736

$ws.Range("H77").Value = 1843.9375
$ws.Range("J77").Value = 2988
$ws.Range("L77").Value = 14940
$ws.Range("N77").Value = -23676

$ws.Range("H136").Value = 2966.2258
$ws.Range("I136").Value = 2423.353
$ws.Range("J136").Value = 3625.4285
$ws.Range("K136").Value = 7270.059
$ws.Range("L136").Value = 10876.2855
$ws.Range("M136").Value = -4720.059
$ws.Range("N136").Value = -15976.2855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 303
$ws.Range("I7").Value = 303
$ws.Range("K7").Value = 303
$ws.Range("M7").Value = -190

$ws.Range("H8").Value = 6337.6665
$ws.Range("J8").Value = 5005
$ws.Range("L8").Value = 5005
$ws.Range("N8").Value = -5285

$ws.Range("H75").Value = 21428.857
$ws.Range("I75").Value = 3089.3333
$ws.Range("J75").Value = 35183.5
$ws.Range("K75").Value = 3089.3333
$ws.Range("L75").Value = 35183.5
$ws.Range("M75").Value = -2153.3333
$ws.Range("N75").Value = -37055.5

$ws.Range("H78").Value = 21428.857
$ws.Range("I78").Value = 3089.3333
$ws.Range("J78").Value = 35183.5
$ws.Range("K78").Value = 9267.999899999999
$ws.Range("L78").Value = 105550.5
$ws.Range("M78").Value = -4587.999899999999
$ws.Range("N78").Value = -114910.5

$ws.Range("H99").Value = 1306.45
$ws.Range("I99").Value = 1252.4166
$ws.Range("J99").Value = 1387.5
$ws.Range("K99").Value = 1252.4166
$ws.Range("L99").Value = 1387.5
$ws.Range("M99").Value = 245.5834
$ws.Range("N99").Value = -4383.5

$ws.Range("H134").Value = 2839.1724
$ws.Range("I134").Value = 2273.5
$ws.Range("J134").Value = 3764.818
$ws.Range("K134").Value = 6820.5
$ws.Range("L134").Value = 11294.454
$ws.Range("M134").Value = -4285.5
$ws.Range("N134").Value = -16364.454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 757.9783
$ws.Range("I113").Value = 492.8125
$ws.Range("J113").Value = 1364.0714
$ws.Range("K113").Value = 1478.4375
$ws.Range("L113").Value = 4092.2142
$ws.Range("M113").Value = 691.5625
$ws.Range("N113").Value = -8432.2142

$ws.Range("H132").Value = 2018.3226
$ws.Range("I132").Value = 1786.2727
$ws.Range("J132").Value = 2145.95
$ws.Range("K132").Value = 16076.4543
$ws.Range("L132").Value = 19313.55
$ws.Range("M132").Value = -13546.4543
$ws.Range("N132").Value = -24373.55

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 550
$ws.Range("I13").Value = 550
$ws.Range("K13").Value = 550
$ws.Range("M13").Value = -411

$ws.Range("H132").Value = 2680.4333
$ws.Range("I132").Value = 2388.7896
$ws.Range("J132").Value = 3184.182
$ws.Range("K132").Value = 7166.3688
$ws.Range("L132").Value = 9552.545999999998
$ws.Range("M132").Value = -4636.3688
$ws.Range("N132").Value = -14612.546

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()

$ws.Range("H46").Value = 771.4286
$ws.Range("I46").Value = 625
$ws.Range("J46").Value = 830
$ws.Range("K46").Value = 625
$ws.Range("L46").Value = 830
$ws.Range("M46").Value = -437
$ws.Range("N46").Value = -1206

$ws.Range("H94").Value = 100000
$ws.Range("J94").Value = 100000
$ws.Range("L94").Value = 100000
$ws.Range("N94").Value = -101352

$ws.Range("H134").Value = 39457.145
$ws.Range("J134").Value = 39457.145
$ws.Range("L134").Value = 39457.145
$ws.Range("N134").Value = -49597.145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1685000.9
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 1685000.9
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 1685000.9
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -1685336.9

$ws.Range("H122").Value = 2253.3157
$ws.Range("I122").Value = 1138
$ws.Range("J122").Value = 3064.4546
$ws.Range("K122").Value = 3414
$ws.Range("L122").Value = 9193.363799999999
$ws.Range("M122").Value = -964
$ws.Range("N122").Value = -14093.3638

$ws.Range("H126").Value = 2600.4
$ws.Range("I126").Value = 2560.8
$ws.Range("J126").Value = 2640
$ws.Range("K126").Value = 7682.400000000001
$ws.Range("L126").Value = 7920
$ws.Range("M126").Value = -5212.400000000001
$ws.Range("N126").Value = -12860
